# Diploma Supplement template fix: delete wrong/duplicate sentence.
#
# Removes the trailing " - <preliminary/previous document...>" clause that
# was mistakenly appended to two label paragraphs:
#   UA: "Інформація про визнання іноземних документів про освіту - попередній
#        документ про освіту, виданий в Україні"
#     -> "Інформація про визнання іноземних документів про освіту"
#   EN: "Information on the recognition of international education documents
#        - preliminary document on education issued in Ukraine"
#     -> "Information on the recognition of international education documents"

$d = $word.ActiveDocument

# --- Ukrainian sentence -------------------------------------------------
$d.Content.Find.Execute(
    "освіту - попередній документ про освіту, виданий в Україні ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "освіту", 2)

# --- English sentence ----------------------------------------------------
$d.Content.Find.Execute(
    "education documents - preliminary document on education issued in Ukraine",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "education documents", 2)
